$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 23 de Marzo de 2020 a las 08:16'

$data = @(
    @('China', 81093, 39, 72703, 5120, 1749, 9, 3270),
    @('Italia', 59138, 0, 7024, 46638, 3000, 0, 5476),
    @('Estados Unidos', 35060, 1514, 178, 34425, 795, 38, 457),
    @('España', 28768, 0, 2575, 24421, 1785, 0, 1772),
    @('Alemania', 24873, 0, 266, 24513, 23, 0, 94),
    @('Iran', 21638, 0, 7913, 12040, 0, 0, 1685),
    @('Francia', 16018, 0, 2200, 13144, 1746, 0, 674),
    @('Corea del Sur', 8961, 64, 3166, 5684, 59, 7, 111),
    @('Suiza', 7474, 0, 131, 7245, 141, 0, 98),
    @('Reino Unido', 5683, 0, 135, 5267, 20, 0, 281),
    @('Paises Bajos', 4204, 0, 2, 4023, 354, 0, 179),
    @('Austria', 3582, 0, 9, 3557, 15, 0, 16),
    @('Belgica', 3401, 0, 340, 2986, 288, 0, 75),
    @('Noruega', 2385, 0, 6, 2372, 28, 0, 7),
    @('Suecia', 1934, 0, 16, 1897, 76, 0, 21),
    @('Australia', 1709, 100, 88, 1614, 11, 0, 7),
    @('Portugal', 1600, 0, 5, 1581, 26, 0, 14),
    @('Brasil', 1546, 0, 2, 1519, 18, 0, 25),
    @('Canada', 1470, 0, 14, 1436, 1, 0, 20),
    @('Dinamarca', 1395, 0, 1, 1381, 42, 0, 13),
    @('Malasia', 1306, 0, 139, 1157, 26, 0, 10),
    @('Israel', 1238, 167, 37, 1200, 24, 0, 1),
    @('Turquia', 1236, 0, 0, 1206, 0, 0, 30),
    @('Chequia', 1120, 0, 6, 1113, 19, 0, 1),
    @('Japon', 1101, 0, 235, 825, 49, 0, 41),
    @('Irlanda', 906, 0, 5, 897, 29, 0, 4),
    @('Pakistan', 799, 23, 13, 780, 0, 1, 6),
    @('Luxemburgo', 798, 0, 6, 784, 3, 0, 8),
    @('Ecuador', 789, 0, 3, 772, 2, 0, 14),
    @('Tailandia', 721, 122, 52, 668, 7, 0, 1),
    @('Crucero', 712, 0, 567, 137, 15, 0, 8),
    @('Polonia', 634, 0, 13, 614, 3, 0, 7),
    @('Chile', 632, 0, 8, 623, 7, 0, 1),
    @('Finlandia', 626, 0, 10, 615, 12, 0, 1),
    @('Grecia', 624, 0, 19, 590, 18, 0, 15),
    @('Islandia', 568, 0, 5, 562, 1, 0, 1),
    @('Indonesia', 514, 0, 29, 437, 0, 0, 48),
    @('Arabia Saudita', 511, 0, 17, 494, 0, 0, 0),
    @('Catar', 494, 0, 33, 461, 6, 0, 0),
    @('Singapur', 455, 0, 144, 309, 14, 0, 2),
    @('Rusia', 438, 71, 16, 421, 0, 0, 1),
    @('Rumania', 433, 0, 64, 366, 14, 0, 3),
    @('India', 425, 29, 24, 393, 0, 1, 8),
    @('Eslovenia', 414, 0, 0, 412, 12, 0, 2),
    @('Filipinas', 396, 16, 18, 345, 1, 8, 33),
    @('Peru', 363, 0, 1, 357, 5, 0, 5),
    @('Barein', 334, 0, 149, 183, 3, 0, 2),
    @('Egipto', 327, 0, 56, 257, 0, 0, 14),
    @('Estonia', 326, 0, 4, 322, 2, 0, 0),
    @('Hong Kong', 318, 0, 100, 214, 4, 0, 4),
    @('Mexico', 316, 65, 4, 310, 1, 0, 2),
    @('Panama', 313, 0, 1, 309, 7, 0, 3),
    @('Sudafrica', 274, 0, 2, 272, 0, 0, 0),
    @('Argentina', 266, 0, 27, 235, 0, 0, 4),
    @('Croacia', 254, 0, 5, 248, 5, 0, 1),
    @('Libano', 248, 0, 8, 236, 4, 0, 4),
    @('Irak', 233, 0, 57, 156, 0, 0, 20),
    @('Colombia', 231, 0, 3, 226, 0, 0, 2),
    @('Serbia', 222, 0, 2, 218, 4, 0, 2),
    @('Republica Dominicana', 202, 0, 0, 199, 0, 0, 3),
    @('Argelia', 201, 0, 65, 119, 0, 0, 17),
    @('Armenia', 194, 0, 2, 192, 6, 0, 0),
    @('Bulgaria', 190, 3, 3, 184, 3, 0, 3),
    @('Kuwait', 188, 0, 30, 158, 5, 0, 0),
    @('Eslovaquia', 185, 0, 7, 178, 2, 0, 0),
    @('San Marino', 175, 0, 4, 151, 13, 0, 20),
    @('Taiwan', 169, 0, 28, 139, 0, 0, 2),
    @('Hungria', 167, 36, 16, 144, 6, 1, 7),
    @('Uruguay', 158, 0, 0, 158, 2, 0, 0),
    @('Emiratos Arabes Unidos', 153, 0, 38, 113, 2, 0, 2),
    @('Lituania', 143, 0, 1, 141, 1, 0, 1),
    @('Letonia', 139, 0, 1, 138, 0, 0, 0),
    @('Costa Rica', 134, 0, 2, 130, 2, 0, 2),
    @('Bosnia y Herzegovina', 126, 0, 2, 123, 1, 0, 1),
    @('Vietnam', 118, 5, 17, 101, 2, 0, 0),
    @('Republica de Macedonia', 115, 0, 1, 113, 1, 0, 1),
    @('Islas Feroe', 115, 0, 3, 112, 0, 0, 0),
    @('Marruecos', 115, 0, 3, 108, 1, 0, 4),
    @('Principado de Andorra', 113, 0, 1, 111, 2, 0, 1),
    @('Jordania', 112, 0, 1, 111, 0, 0, 0),
    @('Nueva Zelanda', 102, 0, 0, 102, 0, 0, 0),
    @('Republica de Chipre', 95, 0, 3, 91, 3, 0, 1),
    @('Moldavia', 94, 0, 1, 92, 3, 0, 1),
    @('Malta', 90, 0, 2, 88, 1, 0, 0),
    @('Albania', 89, 0, 2, 85, 2, 0, 2),
    @('Brunei', 88, 0, 2, 86, 2, 0, 0),
    @('Camboya', 86, 2, 2, 84, 0, 0, 0),
    @('Sri Lanka', 86, 4, 3, 83, 2, 0, 0),
    @('Venezuela', 77, 7, 15, 62, 2, 0, 0),
    @('Bielorrusia', 76, 0, 15, 61, 0, 0, 0),
    @('Tunez', 75, 0, 1, 71, 7, 0, 3),
    @('Burkina Faso', 75, 0, 5, 66, 0, 0, 4),
    @('Ucrania', 73, 0, 1, 69, 0, 0, 3),
    @('Senegal', 67, 0, 5, 62, 0, 0, 0),
    @('Azerbaiyan', 65, 0, 11, 53, 0, 0, 1),
    @('Reunion', 64, 0, 1, 63, 0, 0, 0),
    @('Kazajistan', 60, 0, 0, 60, 0, 0, 0),
    @('Estado de Palestina', 59, 0, 17, 42, 0, 0, 0),
    @('Guadalupe', 58, 0, 0, 57, 4, 0, 1),
    @('Oman', 55, 0, 17, 38, 0, 0, 0),
    @('Georgia', 54, 0, 3, 51, 1, 0, 0),
    @('Trinidad yTobago', 50, 0, 0, 50, 0, 0, 0),
    @('Uzbekistan', 46, 3, 0, 46, 0, 0, 0),
    @('Martinica', 44, 0, 0, 43, 7, 0, 1),
    @('Afganistan', 40, 0, 1, 38, 0, 0, 1),
    @('Camerun', 40, 0, 2, 38, 0, 0, 0),
    @('Liechtenstein', 37, 0, 0, 37, 0, 0, 0),
    @('Cuba', 35, 0, 0, 34, 0, 0, 1),
    @('Consejo Danes para los Refugiados', 30, 0, 0, 29, 0, 0, 1),
    @('Nigeria', 30, 0, 2, 28, 0, 0, 0),
    @('Mauricio', 28, 0, 0, 26, 1, 0, 2),
    @('Guam', 27, 0, 0, 26, 0, 0, 1),
    @('Banglades', 27, 0, 3, 22, 0, 0, 2),
    @('Honduras', 26, 0, 0, 26, 0, 0, 0),
    @('Bolivia', 24, 0, 0, 24, 0, 0, 0),
    @('Ghana', 23, 0, 0, 22, 0, 0, 1),
    @('Monaco', 23, 0, 1, 22, 0, 0, 0),
    @('Puerto Rico', 23, 0, 0, 22, 0, 0, 1),
    @('Montenegro', 22, 1, 0, 21, 0, 1, 1),
    @('Paraguay', 22, 0, 0, 21, 1, 0, 1),
    @('Macao', 22, 0, 10, 12, 0, 0, 0),
    @('Ruanda', 19, 0, 0, 19, 0, 0, 0),
    @('Guyana', 19, 0, 0, 18, 0, 0, 1),
    @('Guatemala', 19, 0, 0, 18, 0, 0, 1),
    @('Jamaica', 19, 0, 2, 16, 0, 0, 1),
    @('Polinesia Francesa', 18, 0, 0, 18, 0, 0, 0),
    @('Guayana Francesa', 18, 0, 0, 18, 0, 0, 0),
    @('Togo', 16, 0, 0, 16, 0, 0, 0),
    @('Kenia', 15, 0, 0, 15, 0, 0, 0),
    @('Gibraltar', 15, 0, 2, 13, 0, 0, 0),
    @('Kirguistan', 14, 0, 0, 14, 0, 0, 0),
    @('Barbados', 14, 0, 0, 14, 0, 0, 0),
    @('Costa de Marfil', 14, 0, 1, 13, 0, 0, 0),
    @('Maldivas', 13, 0, 3, 10, 0, 0, 0),
    @('Madagascar', 12, 9, 0, 12, 0, 0, 0),
    @('Tanzania', 12, 0, 0, 12, 0, 0, 0),
    @('Etiopia', 11, 0, 0, 11, 0, 0, 0),
    @('Mayotte', 11, 0, 0, 11, 0, 0, 0),
    @('Mongolia', 10, 0, 0, 10, 0, 0, 0),
    @('Aruba', 9, 0, 1, 8, 0, 0, 0),
    @('Seychelles', 7, 0, 0, 7, 0, 0, 0),
    @('Bermudas', 6, 0, 0, 6, 0, 0, 0),
    @('Guinea Ecuatorial', 6, 0, 0, 6, 0, 0, 0),
    @('Islas Virgenes de los Estados Unidos', 6, 0, 0, 6, 0, 0, 0),
    @('Surinam', 5, 0, 0, 5, 0, 0, 0),
    @('San Martin (Parte Francesa)', 5, 0, 0, 5, 0, 0, 0),
    @('Isla de Man', 5, 0, 0, 5, 0, 0, 0),
    @('Gabon', 5, 0, 0, 4, 0, 0, 1),
    @('Groenlandia', 4, 0, 0, 4, 0, 0, 0),
    @('Nueva Caledonia', 4, 0, 0, 4, 0, 0, 0),
    @('Suazilandia', 4, 0, 0, 4, 0, 0, 0),
    @('Bahamas', 4, 0, 0, 4, 0, 0, 0),
    @('El Salvador', 3, 0, 0, 3, 0, 0, 0),
    @('Republica de Africa Central', 3, 0, 0, 3, 0, 0, 0),
    @('Zambia', 3, 0, 0, 3, 0, 0, 0),
    @('Liberia', 3, 0, 0, 3, 0, 0, 0),
    @('Congo', 3, 0, 0, 3, 0, 0, 0),
    @('Fiyi', 3, 1, 0, 3, 0, 0, 0),
    @('Zimbabue', 3, 0, 0, 3, 0, 0, 0),
    @('Namibia', 3, 0, 0, 3, 0, 0, 0),
    @('Cabo Verde', 3, 0, 0, 3, 0, 0, 0),
    @('San Bartolome', 3, 0, 0, 3, 0, 0, 0),
    @('Islas Caimanes', 3, 0, 0, 2, 0, 0, 1),
    @('Curazao', 3, 0, 0, 2, 0, 0, 1),
    @('Nicaragua', 2, 0, 0, 2, 0, 0, 0),
    @('Mauritania', 2, 0, 0, 2, 0, 0, 0),
    @('Niger', 2, 0, 0, 2, 0, 0, 0),
    @('Guinea', 2, 0, 0, 2, 0, 0, 0),
    @('Haiti', 2, 0, 0, 2, 0, 0, 0),
    @('Angola', 2, 0, 0, 2, 0, 0, 0),
    @('Santa Lucia', 2, 0, 0, 2, 0, 0, 0),
    @('Butan', 2, 0, 0, 2, 0, 0, 0),
    @('Benin', 2, 0, 0, 2, 0, 0, 0),
    @('Sudan', 2, 0, 0, 1, 0, 0, 1),
    @('Nepal', 2, 1, 1, 1, 0, 0, 0),
    @('Antigua y Barbuda', 1, 0, 0, 1, 0, 0, 0),
    @('Montserrat', 1, 0, 0, 1, 0, 0, 0),
    @('San Vicente y las Granadinas', 1, 0, 0, 1, 0, 0, 0),
    @('Republica de Yibuti', 1, 0, 0, 1, 0, 0, 0),
    @('Granada', 1, 0, 0, 1, 0, 0, 0),
    @('Papua Nueva Guinea', 1, 0, 0, 1, 0, 0, 0),
    @('Timor Oriental', 1, 0, 0, 1, 0, 0, 0),
    @('San Martin (Parte Holandesa)', 1, 0, 0, 1, 0, 0, 0),
    @('Siria', 1, 0, 0, 1, 0, 0, 0),
    @('Dominica', 1, 0, 0, 1, 0, 0, 0),
    @('Eritrea', 1, 0, 0, 1, 0, 0, 0),
    @('Republica del Chad', 1, 0, 0, 1, 0, 0, 0),
    @('Uganda', 1, 0, 0, 1, 0, 0, 0),
    @('Mozambique', 1, 0, 0, 1, 0, 0, 0),
    @('Santa Sede', 1, 0, 0, 1, 0, 0, 0),
    @('Gambia', 1, 0, 0, 1, 0, 0, 0),
    @('Somalia', 1, 0, 0, 1, 0, 0, 0)
)

$r = 4
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

